# "main loop error fixed" - the Won/Lost/Points columns for both groups
# were being filled from the wrong loop iteration; re-point each row at
# its correct values and rename the sheet back to the default "Sheet1".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

# Group A: Won (C), Lost (D), Points (E)
# Group B: Won_B (I), Lost_B (J), Points_B (K)
$values = @{
    2 = @{ C = 2; D = 1; E = 4; I = 2; J = 1; K = 4 }
    3 = @{ C = 1; D = 2; E = 2; I = 0; J = 3; K = 0 }
    4 = @{ C = 0; D = 3; E = 0; I = 3; J = 0; K = 6 }
    5 = @{ C = 3; D = 0; E = 6; I = 1; J = 2; K = 2 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("C$row").Value = $rowVals.C
    $ws.Range("D$row").Value = $rowVals.D
    $ws.Range("E$row").Value = $rowVals.E
    $ws.Range("I$row").Value = $rowVals.I
    $ws.Range("J$row").Value = $rowVals.J
    $ws.Range("K$row").Value = $rowVals.K
}
